# feat: add 2022-Q1 data
#
# The workbook currently has 4 sheets: 2020-Q4, 2021-Q1, 2021-Q2, 总计.
# We need to:
#   1. Insert a new sheet "2022-Q1" (holding the fund-by-fund breakdown for
#      that quarter) right before the "总计" (totals) sheet.
#   2. Prepend a new summary row for "2022-Q1" to the "总计" sheet.
#
# Because this COM host assigns each new sheet's internal sheetId as
# max(existing sheetIds)+1, and we want the new "2022-Q1" sheet to reuse
# sheetId 4 (with "总计" becoming sheetId 5, matching a natural "insert
# before, then re-append" edit), we delete "总计" first, recreate
# "2022-Q1" then "总计" in the right order, and refill 总计's data
# (original 3 rows + the new 2022-Q1 row on top).

$wb = $excel.ActiveWorkbook

# Reference cells (on the untouched "2021-Q2" sheet) whose formatting we
# reuse for the new sheets, so the new cells land on the very same style
# entries the rest of the workbook already uses instead of Excel minting
# fresh (but visually-equivalent) ones:
#   - headerFmt / indexFmt: bold, thin-bordered, centered+top (style used
#     by the header row and the running-index column A)
#   - plainFmt: completely unstyled data cell
$fmtSheet = $wb.Worksheets.Item("2021-Q2")
$headerFmt = $fmtSheet.Range("B1")
$plainFmt = $fmtSheet.Range("C2")

function Copy-Format($srcRange, $dstRange) {
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)   # xlPasteFormats
}

function Set-HeaderCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $text
    Copy-Format $headerFmt $c
    return $c
}

function Set-IndexCell($ws, $row, $col, $n) {
    # Column-A running index (0,1,2,...) styled like the header row.
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $n
    Copy-Format $headerFmt $c
    return $c
}

function Set-TextCell($ws, $row, $col, $text) {
    # Force a cell to store its value as text even when the text looks
    # like a number (e.g. fund codes "910022", figures like "22.15"),
    # then strip the format back down to the workbook's plain/unstyled
    # look (no bold, no border, no explicit number format).
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    Copy-Format $plainFmt $c
    return $c
}

# ---- 1. remove the old "总计" sheet, remembering nothing (data is ------
#        reconstructed from the known original rows below) ---------------

$wb.Worksheets.Item("总计").Delete()

# ---- 2. re-create "2022-Q1" then "总计" at the end, in that order ------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsQ1 = $wb.Worksheets.Add($null, $lastSheet)
$wsQ1.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# ---- 3. fill in "2022-Q1" fund breakdown --------------------------------

Set-HeaderCell $wsQ1 1 2 "基金代码"    | Out-Null
Set-HeaderCell $wsQ1 1 3 "基金名称"    | Out-Null
Set-HeaderCell $wsQ1 1 4 "基金规模"    | Out-Null
Set-HeaderCell $wsQ1 1 5 "股票总仓位"  | Out-Null
Set-HeaderCell $wsQ1 1 6 "仓位占比"    | Out-Null
Set-HeaderCell $wsQ1 1 7 "持有市值(亿元)" | Out-Null
Set-HeaderCell $wsQ1 1 8 "仓位排名"    | Out-Null

$fundRows = @(
    @("910022", "东方红启航三年持有期混合A",       "22.15", "92.82", "2.89", "0.6401", 10),
    @("910028", "东方红内需增长混合型证券投资基金A", "13.44", "92.82", "2.90", "0.3898", 10),
    @("010225", "东方红启航三年持有期混合B",         "13.27", "92.82", "2.89", "0.3835", 10),
    @("012243", "东方红内需增长混合型证券投资基金C", "4.61",  "92.82", "2.90", "0.1337", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $idx = $r - 2
    Set-IndexCell $wsQ1 $r 1 $idx | Out-Null

    Set-TextCell $wsQ1 $r 2 $row[0] | Out-Null

    $cC = $wsQ1.Cells.Item($r, 3)
    $cC.Value = $row[1]

    Set-TextCell $wsQ1 $r 4 $row[2] | Out-Null
    Set-TextCell $wsQ1 $r 5 $row[3] | Out-Null
    Set-TextCell $wsQ1 $r 6 $row[4] | Out-Null
    Set-TextCell $wsQ1 $r 7 $row[5] | Out-Null

    $wsQ1.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---- 4. fill in "总计" (2022-Q1 on top, then the pre-existing rows) ----

Set-HeaderCell $wsTotal 1 2 "日期"           | Out-Null
Set-HeaderCell $wsTotal 1 3 "持有数量(只)"    | Out-Null
Set-HeaderCell $wsTotal 1 4 "持有市值(亿元)"  | Out-Null

$totalRows = @(
    @("2022-Q1", 4,  1.55),
    @("2021-Q2", 2,  0.05),
    @("2021-Q1", 15, 2.03),
    @("2020-Q4", 11, 0.43)
)

$r = 2
foreach ($row in $totalRows) {
    $idx = $r - 2
    Set-IndexCell $wsTotal $r 1 $idx | Out-Null
    $wsTotal.Cells.Item($r, 2).Value = $row[0]
    $wsTotal.Cells.Item($r, 3).Value = $row[1]
    $wsTotal.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
